$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the 'last updated' timestamp in A1
$ws.Range('A1').Value = 'Datos actualizados a 28 de Marzo de 2020 a las 00:28'

# Re-sort the tied rows (149-205) by updating country names / new-case counts
$ws.Range('A149').Value = 'Republica de Yibuti'
$ws.Range('C149').Value = 1
$ws.Range('A150').Value = 'Guinea Ecuatorial'
$ws.Range('C150').Value = 0
$ws.Range('A159').Value = 'Haiti'
$ws.Range('C159').Value = 0
$ws.Range('A160').Value = 'Birmania'
$ws.Range('C160').Value = 3
$ws.Range('A166').Value = 'Mozambique'
$ws.Range('A168').Value = 'Antigua y Barbuda'
$ws.Range('A169').Value = 'Granada'
$ws.Range('A172').Value = 'Laos'
$ws.Range('A173').Value = 'Eritrea'
$ws.Range('A174').Value = 'San Bartolome'
$ws.Range('C174').Value = 2
$ws.Range('A175').Value = 'Fiyi'
$ws.Range('C175').Value = 0
$ws.Range('A177').Value = 'Montserrat'
$ws.Range('A178').Value = 'Guyana'
$ws.Range('A179').Value = 'Zimbabue'
$ws.Range('C179').Value = 2
$ws.Range('A180').Value = 'Cabo Verde'
$ws.Range('C180').Value = 0
$ws.Range('A181').Value = 'Santa Sede'
$ws.Range('A182').Value = 'Angola'
$ws.Range('A183').Value = 'Congo'
$ws.Range('A186').Value = 'Liberia'
$ws.Range('C186').Value = 0
$ws.Range('A187').Value = 'Republica de Africa Central'
$ws.Range('A188').Value = 'Somalia'
$ws.Range('C188').Value = 1
$ws.Range('A189').Value = 'Butan'
$ws.Range('C189').Value = 1
$ws.Range('A190').Value = 'San Martin (Parte Holandesa)'
$ws.Range('C190').Value = 0
$ws.Range('A191').Value = 'Mauritania'
$ws.Range('A192').Value = 'Gambia'
$ws.Range('A193').Value = 'Sudan'
$ws.Range('D193').Value = 0
$ws.Range('H193').Value = 1
$ws.Range('A194').Value = 'Santa Lucia'
$ws.Range('D194').Value = 1
$ws.Range('H194').Value = 0
$ws.Range('A195').Value = 'Islas Turcas y Caicos'
$ws.Range('A196').Value = 'Anguila'
$ws.Range('A198').Value = 'Islas Virgenes Britanicas'
$ws.Range('A199').Value = 'Guinea-Bisau'
$ws.Range('A200').Value = 'San Cristobal y Nieves'
$ws.Range('A202').Value = 'Papua Nueva Guinea'
$ws.Range('A204').Value = 'San Vicente y las Granadinas'
$ws.Range('A205').Value = 'Timor Oriental'
